$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (header numbering "1." .. "20."): columns H:U ("7." .. "20.") are
# removed, leaving empty numeric cells behind (same style, no content).
$ws.Range("H2:U2").ClearContents()

# Row 3 (Temperature): B3:G3 get new readings, H3:U3 are cleared out.
$ws.Range("B3").Value = 21
$ws.Range("C3").Value = 22
$ws.Range("D3").Value = 22
$ws.Range("E3").Value = 21
$ws.Range("F3").Value = 19
$ws.Range("G3").Value = 20
$ws.Range("H3:U3").ClearContents()

# Row 4 (Hour): B4:G4 get new timestamps, H4:U4 are cleared out.
$ws.Range("B4").Value = "21:04:34"
$ws.Range("C4").Value = "21:04:40"
$ws.Range("D4").Value = "21:04:55"
$ws.Range("E4").Value = "21:05:00"
$ws.Range("F4").Value = "21:05:16"
$ws.Range("G4").Value = "21:05:22"
$ws.Range("H4:U4").ClearContents()

# Row 5 (Date): B5:G5 get a new date, H5:U5 are cleared out.
# These cells must keep storing the date as literal text (not get
# auto-converted to a date serial number), so the new value is typed into a
# helper cell that is explicitly formatted as Text first, then copied in via
# PasteSpecial(values only) so the destination cells retain their original
# style/number-format while receiving the literal string.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"
$helper.Value = "06-10-22"
$helper.Copy()
$ws.Range("B5").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C5").PasteSpecial(-4163)
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").PasteSpecial(-4163)
$ws.Range("F5").PasteSpecial(-4163)
$ws.Range("G5").PasteSpecial(-4163)
$helper.ClearContents()
$helper.NumberFormat = "General"
$helper.Clear()
$ws.Range("H5:U5").ClearContents()

# Row 8 (Nominal): H8:U8 cleared (B8:G8 unchanged).
$ws.Range("H8:U8").ClearContents()

# Row 9 (Lower limit): H9:U9 cleared (B9:G9 unchanged).
$ws.Range("H9:U9").ClearContents()

# Row 10 (Upper limit): H10:U10 cleared (B10:G10 unchanged).
$ws.Range("H10:U10").ClearContents()
